$d = $word.ActiveDocument

function Replace-ParagraphXml([string]$findText, [string]$innerXml) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $findText"
    }
    $para = $rng.Paragraphs(1)
    $prng = $para.Range

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'

    [void]$prng.InsertXML($pkg)
}

# 1) "Prajwala T R" -> "Prajwala" + proofErr spellStart/spellEnd + " T R"
$para1 = '<w:p w14:paraId="427CD8D2" w14:textId="77777777" w:rsidR="008A4B94" w:rsidRDefault="008A4B94" w:rsidP="005A65D8">' +
    '<w:pPr>' +
        '<w:pStyle w:val="NormalWeb"/>' +
        '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
        '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>' +
        '<w:textAlignment w:val="baseline"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="222222"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="222222"/>' +
        '</w:rPr>' +
        '<w:t>To,</w:t>' +
    '</w:r>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="222222"/>' +
        '</w:rPr>' +
        '<w:br/>' +
    '</w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r w:rsidRPr="00320699">' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="222222"/>' +
        '</w:rPr>' +
        '<w:t>Prajwala</w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidRPr="00320699">' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="222222"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve"> T R</w:t>' +
    '</w:r>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="222222"/>' +
        '</w:rPr>' +
        '<w:t>,</w:t>' +
    '</w:r>' +
    '</w:p>'
Replace-ParagraphXml "Prajwala T R" $para1

# 2) "EC Campus," -> "EC Campus" + "."
$para2 = '<w:p w14:paraId="14E97045" w14:textId="06532954" w:rsidR="008A4B94" w:rsidRDefault="008A4B94" w:rsidP="00ED53E7">' +
    '<w:pPr>' +
        '<w:pStyle w:val="NormalWeb"/>' +
        '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
        '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>' +
        '<w:jc w:val="both"/>' +
        '<w:textAlignment w:val="baseline"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="222222"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="222222"/>' +
        '</w:rPr>' +
        '<w:t>EC Campus</w:t>' +
    '</w:r>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="222222"/>' +
        '</w:rPr>' +
        '<w:t>.</w:t>' +
    '</w:r>' +
    '</w:p>'
Replace-ParagraphXml "EC Campus," $para2

# 3) "Guram Balaji" -> "Guram" + proofErr spellStart/spellEnd + " Balaji"
$para3 = '<w:p w14:paraId="5B1393B2" w14:textId="2776F1F1" w:rsidR="005A65D8" w:rsidRDefault="00881C24" w:rsidP="00881C24">' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' +
        '<w:t>Guram</w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' +
        '<w:t xml:space="preserve"> Balaji</w:t>' +
    '</w:r>' +
    '<w:r w:rsidR="005A65D8">' +
        '<w:t>,</w:t>' +
    '</w:r>' +
    '</w:p>'
Replace-ParagraphXml "Guram Balaji" $para3

Write-Host "Done"
